$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Shuffled_Rand_removed")

# New column AG: header + fill color (red) for the data block AG1:AG44
$ws.Range("AG1:AG44").Interior.Color = 255

# Header
$ws.Range("AG1").Value = "Test_by_NN_Model_Ver1.1"

# Data values (rows 2-44)
$ws.Range("AG2").Value = 68.751419999999996
$ws.Range("AG3").Value = 55.454532999999998
$ws.Range("AG4").Value = -5.3729706000000004
$ws.Range("AG5").Value = 22.294035000000001
$ws.Range("AG6").Value = 33.027996000000002
$ws.Range("AG7").Value = 34.922671999999999
$ws.Range("AG8").Value = 126.76423
$ws.Range("AG9").Value = 30.668865
$ws.Range("AG10").Value = 112.17992
$ws.Range("AG11").Value = 27.497883000000002
$ws.Range("AG12").Value = 13.376815000000001
$ws.Range("AG13").Value = 100.48035
$ws.Range("AG14").Value = 122.66893
$ws.Range("AG15").Value = 88.287270000000007
$ws.Range("AG16").Value = 24.917252999999999
$ws.Range("AG17").Value = 79.326099999999997
$ws.Range("AG18").Value = 106.89429
$ws.Range("AG19").Value = 64.694360000000003
$ws.Range("AG20").Value = 30.91028
$ws.Range("AG21").Value = 45.266629999999999
$ws.Range("AG22").Value = 25.980103
$ws.Range("AG23").Value = 61.799045999999997
$ws.Range("AG24").Value = 59.996014000000002
$ws.Range("AG25").Value = 92.902169999999998
$ws.Range("AG26").Value = 4.6627280000000004
$ws.Range("AG27").Value = 75.546239999999997
$ws.Range("AG28").Value = 71.204864999999998
$ws.Range("AG29").Value = 4.8061290000000003
$ws.Range("AG30").Value = 132.30789999999999
$ws.Range("AG31").Value = 31.062002
$ws.Range("AG32").Value = 127.54550999999999
$ws.Range("AG33").Value = 101.22172500000001
$ws.Range("AG34").Value = 43.221989999999998
$ws.Range("AG35").Value = 50.949820000000003
$ws.Range("AG36").Value = 42.059852999999997
$ws.Range("AG37").Value = 6.0401315999999996
$ws.Range("AG38").Value = 40.512529999999998
$ws.Range("AG39").Value = -14.991389
$ws.Range("AG40").Value = 75.581695999999994
$ws.Range("AG41").Value = 54.924610000000001
$ws.Range("AG42").Value = 14.715389999999999
$ws.Range("AG43").Value = 60.149340000000002
$ws.Range("AG44").Value = 13.056656

# Selection as left by the editor
$ws.Activate()
$ws.Range("R19").Select()
